# Generate Report for Handoff
# The c0d98220-... file has finished translation and is ready to hand off.
# Update its status/priority/timestamps on the Overview, zh-cn and de-de
# sheets (row 3 in each table is the c0d98220 file).

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-25 06:13:16"

# --- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-25 06:13:11"

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-25 06:13:16"

# --- Widen the Status columns to fit the new, longer text ------------
# "Ready for handoff" is wider than "In Translation" so Excel's
# column-autofit grows these columns when the sheet is next displayed.
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 16.33
$ws.Columns.Item(6).ColumnWidth = 16.33

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(3).ColumnWidth = 16.33

$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(3).ColumnWidth = 16.33
